$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the value of A5 (tiny precision change)
$ws.Range("A5").Value = 45875.2085046412

# Add new row 6 data
$ws.Range("A6").Value = 45875.25017235359
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = 23
$ws.Range("D6").Value = 13.72
$ws.Range("E6").Value = 92.84
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = "-"
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "06:00:14"
